$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number, Week-covering dates) ---
# "Volume 30   Number  18" -> "...19" (chars 21-22 are "18")
$ws.Range("A8").Characters(21, 2).Text = "19"
# "Report Covering the Week  5/1/2023  Through  5/7/2023"
# -> "...5/8/2023  Through  5/14/2023" (do fixed-length date first, chars unaffected by the other replace)
$ws.Range("C9").Characters(27, 8).Text = "5/8/2023"
$ws.Range("C9").Characters(46, 8).Text = "5/14/2023"

# --- Row 15 ---
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -20
$ws.Range("N15").Value = -52.941176470588

# --- Row 16 ---
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -19.047619047619
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = -35.353535353535
$ws.Range("L16").Value = 30.612244897959
$ws.Range("M16").Value = 6.666666666666
$ws.Range("N16").Value = -78.737541528239

# --- Row 17 ---
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 600
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 136.363636363636
$ws.Range("I17").Value = 89
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = 32.835820895522
$ws.Range("L17").Value = 53.448275862069
$ws.Range("M17").Value = 61.818181818181
$ws.Range("N17").Value = -54.123711340206

# --- Row 18 ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 126
$ws.Range("K18").Value = -20.63492063492
$ws.Range("L18").Value = -18.032786885245
$ws.Range("M18").Value = 8.695652173913
$ws.Range("N18").Value = -66.101694915254

# --- Row 19 ---
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 46.666666666666
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 60.714285714285
$ws.Range("I19").Value = 368
$ws.Range("J19").Value = 330
$ws.Range("K19").Value = 11.515151515151
$ws.Range("L19").Value = 84.924623115577
$ws.Range("M19").Value = 32.851985559566
$ws.Range("N19").Value = -29.230769230769

# --- Row 20 ---
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 1
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = 0
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = 7.142857142857
$ws.Range("L20").Value = -34.782608695652
$ws.Range("M20").Value = -6.25
$ws.Range("N20").Value = -91.891891891891

# --- Row 21 ---
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 36.666666666666
$ws.Range("F21").Value = 159
$ws.Range("G21").Value = 120
$ws.Range("H21").Value = 32.5
$ws.Range("I21").Value = 644
$ws.Range("J21").Value = 649
$ws.Range("K21").Value = -0.770416024653
$ws.Range("L21").Value = 39.393939393939
$ws.Range("M21").Value = 27.524752475247
$ws.Range("N21").Value = -57.519788918205

# --- Row 22 ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -14.285714285714
$ws.Range("L22").Value = 100

# --- Row 23 ---
$ws.Range("C23").Value = 3
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 333.333333333333
$ws.Range("I23").Value = 52
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -29.729729729729
$ws.Range("M23").Value = 20.930232558139

# --- Row 24 ---
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -11.764705882352
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 171
$ws.Range("H24").Value = -28.654970760233
$ws.Range("I24").Value = 536
$ws.Range("J24").Value = 786
$ws.Range("K24").Value = -31.806615776081
$ws.Range("L24").Value = 72.903225806451
$ws.Range("M24").Value = -4.964539007092

# --- Row 25 ---
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 171.428571428571
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 16.216216216216
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 166
$ws.Range("K25").Value = 7.22891566265
$ws.Range("L25").Value = 69.523809523809
$ws.Range("M25").Value = 12.658227848101

# --- Row 26 ---
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -40
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -16.666666666666

# --- Row 27 ---
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -36.363636363636
$ws.Range("L27").Value = 23.529411764705

# --- Row 28 ---
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = -84.615384615384

# --- Row 29 ---
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 1
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -77.777777777777

